# complex_validate_test.xlsx — remove the 'default' column (column R) from the
# three "prompt" sheets (survey, section1, section2). The 'default' header
# (shared string "default") sat between the remaining prompt columns and the
# trailing "hideInContents" / "validation_tags" columns, so deleting it shifts
# everything after it one column to the left.
$wb = $excel.ActiveWorkbook

$sheetNames = @("survey", "section1", "section2")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    # Select the whole 'default' column first (mirrors a user right-click ->
    # Delete on the column header), then remove it.
    $ws.Range("R:R").Select()
    $ws.Range("R1").EntireColumn.Delete()
}

# Leave the selection/cursor on the last-touched sheet (section2) at P10, and
# make that sheet the active tab, matching where the edit was finished.
$ws3 = $wb.Worksheets.Item("section2")
$ws3.Range("P10").Select()
$ws3.Activate()
